# Croatia model input variables: update recalibrated trajectory values
# (gdp_mmm_usd, climate_change_factor_gnrl_hydropower_availability,
#  elasticity_gnrl_rate_occupancy_to_gdppc, frac_gnrl_eating_red_meat)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("strategy_id-0")
$ws2 = $wb.Worksheets.Item("strategy_id-6008")

function Set-RowValues($ws, $row, $startCol, $values) {
    for ($i = 0; $i -lt $values.Count; $i++) {
        $ws.Cells.Item($row, $startCol + $i).Value = $values[$i]
    }
}

# strategy_id-0 : row 2 - gdp_mmm_usd (J2:AS2)
$row2Values = @(45.65220107373895,46.41036762541388,47.16853417708881,47.92670072876376,48.68486728043868,49.44303383211361,55.89097083028935,59.81855567550635,61.64952464147393,63.95369176876388,65.80045132624767,67.58328810903468,69.49299557107426,71.43436099947894,73.22895007874062,75.02130488319007,76.86666342502818,78.53728796439511,79.9829569288255,81.4234330266046,82.99949792738515,84.61742092185833,86.18486613322112,87.67209125245746,89.13240235490146,90.57669295724192,92.00185579690366,93.37649006996322,94.69880878461068,95.97264182058078,97.22973575281335,98.46577926210448,99.66362480501522,100.8345592737375,101.9715451180892,103.0789704758992)
Set-RowValues $ws1 2 10 $row2Values

# strategy_id-0 : row 4 - climate_change_factor_gnrl_hydropower_availability (J4:AS4)
$row4Values = @(1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1)
Set-RowValues $ws1 4 10 $row4Values

# strategy_id-0 : row 5 - elasticity_gnrl_rate_occupancy_to_gdppc (J5:AS5)
$row5Values = @(-0.1,-0.1,-0.1,-0.1,-0.1,-0.1,-0.1,-0.1,-0.1,-0.1,-0.1,-0.1,-0.1,-0.1,-0.1,-0.1,-0.1,-0.1,-0.1,-0.1,-0.1,-0.1,-0.1,-0.1,-0.1,-0.1,-0.1,-0.1,-0.1,-0.1,-0.1,-0.1,-0.1,-0.1,-0.1,-0.1)
Set-RowValues $ws1 5 10 $row5Values

# strategy_id-0 : row 6 - frac_gnrl_eating_red_meat (J6:AS6, J6 unchanged at 1)
$row6Values = @(1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1)
Set-RowValues $ws1 6 10 $row6Values

# strategy_id-6008 : row 2 - frac_gnrl_eating_red_meat (J2:AS2, J2 unchanged at 1)
$row2bValues = @(1,1,1,1,1,1,1,1,1,1,1,1,1,0.9782608695652174,0.9565217391304348,0.9347826086956521,0.9130434782608696,0.891304347826087,0.8695652173913043,0.8478260869565217,0.8260869565217391,0.8043478260869565,0.7826086956521739,0.7608695652173914,0.7391304347826086,0.7173913043478262,0.6956521739130435,0.6739130434782609,0.6521739130434783,0.6304347826086957,0.6086956521739131,0.5869565217391304,0.5652173913043479,0.5434782608695652,0.5217391304347826,0.5)
Set-RowValues $ws2 2 10 $row2bValues

